$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "SAN DIEGO AREA TOTALS" label from B2 to A2, and set B2 to "Totals"
$label = $ws.Range("B2").Text
$ws.Range("A2").Value = $label
$ws.Range("B2").Value = "Totals"

# Update the active cell selection to B3
$ws.Range("B3").Select()
